$d = $word.ActiveDocument

$pairs = @(
    @("31×29=899", "99×62=6138"),
    @("35×68=2380", "75×97=7275"),
    @("77×12=924", "46×15=690"),
    @("77×33=2541", "36×47=1692"),
    @("25×12=300", "72×59=4248"),
    @("96×52=4992", "30×72=2160"),
    @("15×51=765", "13×30=390"),
    @("85×53=4505", "73×53=3869"),
    @("34×35=1190", "86×29=2494"),
    @("74×29=2146", "92×72=6624"),
    @("70×62=4340", "19×26=494"),
    @("24×80=1920", "29×26=754"),
    @("54×51=2754", "17×31=527"),
    @("43×85=3655", "55×34=1870"),
    @("19×42=798", "75×31=2325"),
    @("65×89=5785", "68×18=1224"),
    @("65×34=2210", "68×56=3808"),
    @("87×49=4263", "52×72=3744"),
    @("46×54=2484", "57×45=2565"),
    @("66×82=5412", "75×16=1200"),
    @("53×15=795", "22×11=242"),
    @("43×15=645", "59×63=3717"),
    @("44×15=660", "33×49=1617"),
    @("20×66=1320", "71×45=3195"),
    @("27×42=1134", "21×73=1533")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
